$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The updated analysis now covers 3 sending clusters x 3 target clusters (9 data rows)
# instead of the previous 2 sending clusters x 4 target clusters (8 data rows).
# Clear the old data block (rows 2-9) before writing the refreshed TPM-derived table.
$ws.Range("A2:T9").ClearContents()

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efna5"
$ws.Cells.Item(2, 3).Value = "Epha3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.114918
$ws.Cells.Item(2, 8).Value = 0.344754
$ws.Cells.Item(2, 9).Value = 0.04640425382421802
$ws.Cells.Item(2, 10).Value = 0.04640425382421801
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.114581
$ws.Cells.Item(2, 14).Value = 0.343743
$ws.Cells.Item(2, 15).Value = 0.007635610029470834
$ws.Cells.Item(2, 16).Value = 0.007635610029470834
$ws.Cells.Item(2, 17).Value = 0.013167419358
$ws.Cells.Item(2, 18).Value = 0.118506774222
$ws.Cells.Item(2, 19).Value = 0.0003543247859103094
$ws.Cells.Item(2, 20).Value = 0.0003543247859103094

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efna5"
$ws.Cells.Item(3, 3).Value = "Epha3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.114918
$ws.Cells.Item(3, 8).Value = 0.344754
$ws.Cells.Item(3, 9).Value = 0.04640425382421802
$ws.Cells.Item(3, 10).Value = 0.04640425382421801
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 14.85444966666667
$ws.Cells.Item(3, 14).Value = 44.563349
$ws.Cells.Item(3, 15).Value = 0.9898917347297518
$ws.Cells.Item(3, 16).Value = 0.9898917347297518
$ws.Cells.Item(3, 17).Value = 1.707043646794
$ws.Cells.Item(3, 18).Value = 15.363392821146
$ws.Cells.Item(3, 19).Value = 0.04593518731689489
$ws.Cells.Item(3, 20).Value = 0.04593518731689489

# Row 4: ECs -> MuSCs
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efna5"
$ws.Cells.Item(4, 3).Value = "Epha3"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.114918
$ws.Cells.Item(4, 8).Value = 0.344754
$ws.Cells.Item(4, 9).Value = 0.04640425382421802
$ws.Cells.Item(4, 10).Value = 0.04640425382421801
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.037105
$ws.Cells.Item(4, 14).Value = 0.111315
$ws.Cells.Item(4, 15).Value = 0.0024726552407774
$ws.Cells.Item(4, 16).Value = 0.0024726552407774
$ws.Cells.Item(4, 17).Value = 0.00426403239
$ws.Cells.Item(4, 18).Value = 0.03837629151
$ws.Cells.Item(4, 19).Value = 0.0001147417214128174
$ws.Cells.Item(4, 20).Value = 0.0001147417214128174

# Row 5: FAPs -> ECs
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Efna5"
$ws.Cells.Item(5, 3).Value = "Epha3"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.030023666666667
$ws.Cells.Item(5, 8).Value = 6.090071
$ws.Cells.Item(5, 9).Value = 0.819730011809897
$ws.Cells.Item(5, 10).Value = 0.819730011809897
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.114581
$ws.Cells.Item(5, 14).Value = 0.343743
$ws.Cells.Item(5, 15).Value = 0.007635610029470834
$ws.Cells.Item(5, 16).Value = 0.007635610029470834
$ws.Cells.Item(5, 17).Value = 0.2326021417503333
$ws.Cells.Item(5, 18).Value = 2.093419275753
$ws.Cells.Item(5, 19).Value = 0.006259138699633895
$ws.Cells.Item(5, 20).Value = 0.006259138699633895

# Row 6: FAPs -> FAPs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Efna5"
$ws.Cells.Item(6, 3).Value = "Epha3"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.030023666666667
$ws.Cells.Item(6, 8).Value = 6.090071
$ws.Cells.Item(6, 9).Value = 0.819730011809897
$ws.Cells.Item(6, 10).Value = 0.819730011809897
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 14.85444966666667
$ws.Cells.Item(6, 14).Value = 44.563349
$ws.Cells.Item(6, 15).Value = 0.9898917347297518
$ws.Cells.Item(6, 16).Value = 0.9898917347297518
$ws.Cells.Item(6, 17).Value = 30.15488437864211
$ws.Cells.Item(6, 18).Value = 271.393959407779
$ws.Cells.Item(6, 19).Value = 0.8114439634005389
$ws.Cells.Item(6, 20).Value = 0.8114439634005389

# Row 7: FAPs -> MuSCs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Efna5"
$ws.Cells.Item(7, 3).Value = "Epha3"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.030023666666667
$ws.Cells.Item(7, 8).Value = 6.090071
$ws.Cells.Item(7, 9).Value = 0.819730011809897
$ws.Cells.Item(7, 10).Value = 0.819730011809897
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.037105
$ws.Cells.Item(7, 14).Value = 0.111315
$ws.Cells.Item(7, 15).Value = 0.0024726552407774
$ws.Cells.Item(7, 16).Value = 0.0024726552407774
$ws.Cells.Item(7, 17).Value = 0.07532402815166667
$ws.Cells.Item(7, 18).Value = 0.6779162533649999
$ws.Cells.Item(7, 19).Value = 0.002026909709724262
$ws.Cells.Item(7, 20).Value = 0.002026909709724262

# Row 8: MuSCs -> ECs
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Efna5"
$ws.Cells.Item(8, 3).Value = "Epha3"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.3315123333333334
$ws.Cells.Item(8, 8).Value = 0.994537
$ws.Cells.Item(8, 9).Value = 0.133865734365885
$ws.Cells.Item(8, 10).Value = 0.133865734365885
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.114581
$ws.Cells.Item(8, 14).Value = 0.343743
$ws.Cells.Item(8, 15).Value = 0.007635610029470834
$ws.Cells.Item(8, 16).Value = 0.007635610029470834
$ws.Cells.Item(8, 17).Value = 0.03798501466566667
$ws.Cells.Item(8, 18).Value = 0.341865131991
$ws.Cells.Item(8, 19).Value = 0.00102214654392663
$ws.Cells.Item(8, 20).Value = 0.00102214654392663

# Row 9: MuSCs -> FAPs
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Efna5"
$ws.Cells.Item(9, 3).Value = "Epha3"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.3315123333333334
$ws.Cells.Item(9, 8).Value = 0.994537
$ws.Cells.Item(9, 9).Value = 0.133865734365885
$ws.Cells.Item(9, 10).Value = 0.133865734365885
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 14.85444966666667
$ws.Cells.Item(9, 14).Value = 44.563349
$ws.Cells.Item(9, 15).Value = 0.9898917347297518
$ws.Cells.Item(9, 16).Value = 0.9898917347297518
$ws.Cells.Item(9, 17).Value = 4.924433269379223
$ws.Cells.Item(9, 18).Value = 44.319899424413
$ws.Cells.Item(9, 19).Value = 0.1325125840123181
$ws.Cells.Item(9, 20).Value = 0.132512584012318

# Row 10: MuSCs -> MuSCs
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Efna5"
$ws.Cells.Item(10, 3).Value = "Epha3"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.3315123333333334
$ws.Cells.Item(10, 8).Value = 0.994537
$ws.Cells.Item(10, 9).Value = 0.133865734365885
$ws.Cells.Item(10, 10).Value = 0.133865734365885
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.037105
$ws.Cells.Item(10, 14).Value = 0.111315
$ws.Cells.Item(10, 15).Value = 0.0024726552407774
$ws.Cells.Item(10, 16).Value = 0.0024726552407774
$ws.Cells.Item(10, 17).Value = 0.01230076512833333
$ws.Cells.Item(10, 18).Value = 0.110706886155
$ws.Cells.Item(10, 19).Value = 0.0003310038096403208
$ws.Cells.Item(10, 20).Value = 0.0003310038096403208

Write-Host "done"